$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell updates derived from the authoritative diff: coin price/volume refresh
# plus a few row re-orderings (rank swaps) on 2024-03-26.
$updates = @(
    @{Row=2; Col="D"; Val="70.007.15"}
    @{Row=3; Col="D"; Val="3.576.22"}
    @{Row=3; Col="E"; Val="  +0.50%  "}
    @{Row=4; Col="E"; Val="  -0.11%  "}
    @{Row=5; Col="D"; Val="'577.49"}
    @{Row=5; Col="E"; Val="  -1.98%  "}
    @{Row=6; Col="D"; Val="'190.59"}
    @{Row=6; Col="E"; Val="  -0.61%  "}
    @{Row=7; Col="D"; Val="'0.632"}
    @{Row=7; Col="E"; Val="  -1.51%  "}
    @{Row=8; Col="D"; Val="3.568.48"}
    @{Row=8; Col="E"; Val="  +0.51%  "}
    @{Row=9; Col="E"; Val="  +0.00%  "}
    @{Row=10; Col="D"; Val="'0.177"}
    @{Row=10; Col="E"; Val="  -2.40%  "}
    @{Row=11; Col="D"; Val="'0.659"}
    @{Row=11; Col="E"; Val="  +0.01%  "}
    @{Row=12; Col="D"; Val="'56.59"}
    @{Row=12; Col="E"; Val="  -2.66%  "}
    @{Row=13; Col="D"; Val="'0.0000299"}
    @{Row=13; Col="E"; Val="  +2.27%  "}
    @{Row=14; Col="D"; Val="'9.80"}
    @{Row=14; Col="E"; Val="  +1.62%  "}
    @{Row=15; Col="D"; Val="4.154.82"}
    @{Row=15; Col="E"; Val="  +0.74%  "}
    @{Row=16; Col="D"; Val="'20.15"}
    @{Row=16; Col="E"; Val="  +4.95%  "}
    @{Row=17; Col="D"; Val="3.574.77"}
    @{Row=17; Col="E"; Val="  +0.28%  "}
    @{Row=18; Col="D"; Val="69.927.66"}
    @{Row=18; Col="E"; Val="  +0.93%  "}
    @{Row=19; Col="D"; Val="'12.49"}
    @{Row=19; Col="E"; Val="  +0.74%  "}
    @{Row=20; Col="D"; Val="'0.121"}
    @{Row=20; Col="E"; Val="  +1.22%  "}
    @{Row=21; Col="D"; Val="'1.04"}
    @{Row=21; Col="E"; Val="  -0.14%  "}
    @{Row=22; Col="D"; Val="'20.12"}
    @{Row=22; Col="E"; Val="  +18.13%  "}
    @{Row=23; Col="D"; Val="'476.81"}
    @{Row=23; Col="E"; Val="  -5.48%  "}
    @{Row=24; Col="D"; Val="'5.09"}
    @{Row=24; Col="E"; Val="  -7.76%  "}
    @{Row=25; Col="D"; Val="'4.35"}
    @{Row=25; Col="E"; Val="  -2.12%  "}
    @{Row=26; Col="D"; Val="'88.66"}
    @{Row=26; Col="E"; Val="  -2.58%  "}
    @{Row=27; Col="D"; Val="'3.06"}
    @{Row=27; Col="E"; Val="  +1.11%  "}
    @{Row=28; Col="D"; Val="'11.16"}
    @{Row=28; Col="E"; Val="  -0.07%  "}
    @{Row=29; Col="D"; Val="'9.28"}
    @{Row=29; Col="E"; Val="  +0.30%  "}
    @{Row=30; Col="D"; Val="'7.76"}
    @{Row=30; Col="E"; Val="  +3.94%  "}
    @{Row=31; Col="D"; Val="'32.09"}
    @{Row=31; Col="E"; Val="  +0.62%  "}
    @{Row=32; Col="D"; Val="'0.121"}
    @{Row=32; Col="E"; Val="  +5.89%  "}
    @{Row=33; Col="B"; Val="Cosmos"}
    @{Row=33; Col="C"; Val="https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"}
    @{Row=33; Col="D"; Val="'12.12"}
    @{Row=33; Col="E"; Val="  +0.04%  "}
    @{Row=34; Col="B"; Val="OKB"}
    @{Row=34; Col="C"; Val="https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"}
    @{Row=34; Col="D"; Val="'65.67"}
    @{Row=34; Col="E"; Val="  +0.63%  "}
    @{Row=35; Col="D"; Val="'597.33"}
    @{Row=35; Col="E"; Val="  -3.07%  "}
    @{Row=36; Col="D"; Val="'39.88"}
    @{Row=36; Col="E"; Val="  +5.66%  "}
    @{Row=37; Col="D"; Val="0.0₃0807"}
    @{Row=37; Col="E"; Val="  -1.94%  "}
    @{Row=38; Col="B"; Val="Dai"}
    @{Row=38; Col="C"; Val="https://coinranking.com/coin/MoTuySvg7+dai-dai"}
    @{Row=38; Col="D"; Val="'1.00"}
    @{Row=38; Col="E"; Val="  +0.02%  "}
    @{Row=39; Col="B"; Val="TheGraph"}
    @{Row=39; Col="C"; Val="https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"}
    @{Row=39; Col="D"; Val="'0.401"}
    @{Row=39; Col="E"; Val="  +1.37%  "}
    @{Row=40; Col="B"; Val="Kaspa"}
    @{Row=40; Col="C"; Val="https://coinranking.com/coin/V8GxkwWow+kaspa-kas"}
    @{Row=40; Col="D"; Val="'0.145"}
    @{Row=40; Col="E"; Val="  -2.09%  "}
    @{Row=41; Col="D"; Val="'2.96"}
    @{Row=41; Col="E"; Val="  +10.01%  "}
    @{Row=42; Col="D"; Val="'3.54"}
    @{Row=42; Col="E"; Val="  -2.67%  "}
    @{Row=43; Col="B"; Val="ThetaToken"}
    @{Row=43; Col="C"; Val="https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"}
    @{Row=43; Col="D"; Val="'3.17"}
    @{Row=43; Col="E"; Val="  +3.84%  "}
    @{Row=44; Col="B"; Val="Maker"}
    @{Row=44; Col="C"; Val="https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"}
    @{Row=44; Col="D"; Val="3.235.87"}
    @{Row=44; Col="E"; Val="  -2.37%  "}
    @{Row=45; Col="B"; Val="dogwifhat"}
    @{Row=45; Col="C"; Val="https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"}
    @{Row=45; Col="D"; Val="'3.12"}
    @{Row=45; Col="E"; Val="  +7.45%  "}
    @{Row=46; Col="D"; Val="'0.0446"}
    @{Row=46; Col="E"; Val="  +1.35%  "}
    @{Row=47; Col="D"; Val="'9.58"}
    @{Row=47; Col="E"; Val="  +5.70%  "}
    @{Row=48; Col="D"; Val="'3.36"}
    @{Row=48; Col="E"; Val="  +2.59%  "}
    @{Row=49; Col="D"; Val="'0.137"}
    @{Row=49; Col="E"; Val="  -0.23%  "}
    @{Row=50; Col="D"; Val="'0.998"}
    @{Row=50; Col="E"; Val="  -0.24%  "}
    @{Row=51; Col="D"; Val="'3.16"}
    @{Row=51; Col="E"; Val="  -2.19%  "}
)

foreach ($u in $updates) {
    $ws.Range($u.Col + $u.Row).Value = $u.Val
}
